# Update cryptocurrency price/volume data
# ("Updated cryptos list on Thu Mar 23 07:12:24 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "324.22", "0.9992").
# A bare Range.Value assignment lets Excel's COM layer auto-coerce those
# strings to floating point numbers, which silently destroys the fixed
# number of decimals / leading zeros the site renders (0.07480 -> 0.0748).
# Prefixing with an apostrophe forces a literal text entry, exactly like
# typing '0.07480 into the Excel UI, so the text round-trips exactly.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.723.09"
$ws.Range("E2").Value = "  -2.04%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.754.61"
$ws.Range("E3").Value = "  -2.69%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'324.22"
$ws.Range("E5").Value = "  -4.68%  "

# Row 6 - USDC
$ws.Range("D6").Value = "'0.9992"

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4313"
$ws.Range("E7").Value = "  -6.05%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3655"
$ws.Range("E8").Value = "  -4.34%  "

# Row 9 - OKB
$ws.Range("D9").Value = "'45.26"
$ws.Range("E9").Value = "  +0.16%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.07474"
$ws.Range("E10").Value = "  -1.74%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "'1.121"
$ws.Range("E11").Value = "  -3.10%  "

# Row 12 - BinanceUSD
$ws.Range("D12").Value = "'0.9993"
$ws.Range("E12").Value = "  -0.28%  "

# Row 13 - Solana
$ws.Range("D13").Value = "'21.62"
$ws.Range("E13").Value = "  -4.25%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'6.165"
$ws.Range("E14").Value = "  -3.38%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'7.256"
$ws.Range("E15").Value = "  -3.98%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "1.752.73"
$ws.Range("E16").Value = "  -2.88%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "'0.00001067"
$ws.Range("E17").Value = "  -2.65%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "'88.36"
$ws.Range("E18").Value = "  +8.14%  "

# Row 19 - TRON
$ws.Range("D19").Value = "'0.06215"
$ws.Range("E19").Value = "  -7.77%  "

# Row 20 - Dai
$ws.Range("D20").Value = "'0.9990"

# Row 21 - Avalanche
$ws.Range("D21").Value = "'17.14"
$ws.Range("E21").Value = "  -2.20%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'6.165"
$ws.Range("E22").Value = "  -4.27%  "

# Row 23 - BitDAO
$ws.Range("D23").Value = "'0.5285"
$ws.Range("E23").Value = "  -5.48%  "

# Row 24 - WrappedBTC
$ws.Range("D24").Value = "27.729.47"
$ws.Range("E24").Value = "  -2.01%  "

# Row 25 - Cosmos
$ws.Range("D25").Value = "'11.66"
$ws.Range("E25").Value = "  -2.25%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "'2.329"
$ws.Range("E26").Value = "  -3.75%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'20.61"
$ws.Range("E27").Value = "  -0.82%  "

# Row 28 - Monero
$ws.Range("D28").Value = "'153.09"
$ws.Range("E28").Value = "  -0.22%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'2.371"
$ws.Range("E29").Value = "  -0.64%  "

# Row 30 - WrappedliquidstakedEther2.0
$ws.Range("D30").Value = "1.950.76"
$ws.Range("E30").Value = "  -3.01%  "

# Row 31 - ImmutableX
$ws.Range("D31").Value = "'1.226"
$ws.Range("E31").Value = "  -2.14%  "

# Row 32 - BitcoinCash
$ws.Range("D32").Value = "'127.46"

# Row 33 - Filecoin
$ws.Range("D33").Value = "'5.733"
$ws.Range("E33").Value = "  -2.50%  "

# Row 34 - Stellar
$ws.Range("D34").Value = "'0.09159"
$ws.Range("E34").Value = "  -4.01%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "'3.653"
$ws.Range("E35").Value = "  -9.42%  "

# Row 36 - Aptos
$ws.Range("D36").Value = "'12.69"
$ws.Range("E36").Value = "  +4.36%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.02314"
$ws.Range("E37").Value = "  -1.98%  "

# Row 38 - Algorand
$ws.Range("D38").Value = "'0.2160"
$ws.Range("E38").Value = "  -6.52%  "

# Row 39 - InternetComputer(DFINITY)
$ws.Range("D39").Value = "'5.118"
$ws.Range("E39").Value = "  -3.46%  "

# Row 40 - TheSandbox
$ws.Range("D40").Value = "'0.6492"
$ws.Range("E40").Value = "  -2.49%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "'0.06105"
$ws.Range("E41").Value = "  -4.06%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "'1.198"
$ws.Range("E42").Value = "  -3.40%  "

# Row 43 - WEMIXTOKEN
$ws.Range("D43").Value = "'1.428"
$ws.Range("E43").Value = "  -3.87%  "

# Row 44 - FraxShare
$ws.Range("D44").Value = "'7.972"
$ws.Range("E44").Value = "  -5.49%  "

# Row 45 - Frax
$ws.Range("D45").Value = "'0.9985"
$ws.Range("E45").Value = "  -0.15%  "

# Row 46 - EnergySwap
$ws.Range("D46").Value = "'13.83"
$ws.Range("E46").Value = "  -2.92%  "

# Row 49 - Quant
$ws.Range("D49").Value = "'126.02"
$ws.Range("E49").Value = "  -3.79%  "

# Row 50 - NEARProtocol
$ws.Range("D50").Value = "'1.980"
$ws.Range("E50").Value = "  -3.25%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "'0.06900"
$ws.Range("E51").Value = "  -3.97%  "

# Row 47 - was PancakeSwap, now Decentraland (coins swapped order with row 48)
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5953"
$ws.Range("E47").Value = "  -3.25%  "

# Row 48 - was Decentraland, now PancakeSwap (coins swapped order with row 47)
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.750"
$ws.Range("E48").Value = "  -3.32%  "

